# Update racial misclassification language to neutral terms.
#
# Three locations change "affecting all Black and Asian-American voters"
# (or variants) to neutral "affecting 50M voters" language:
#   1. PROFESSIONAL SUMMARY paragraph (plain text swap, single run)
#   2. Siege Analytics bullet point (text swap + "50M" becomes its own
#      bold, colored run, matching the formatting already used for the
#      adjacent "23%"/"64%" metric runs in that same paragraph)
#   3. "Impact:" project statement (plain text swap, with "nationwide"
#      appended, single run)

$d = $word.ActiveDocument

# --- 1. Professional summary paragraph --------------------------------
# Stays a single plain-text run -- no new formatting here.
$r1 = $d.Content.Find.Execute(
    "Discovered systematic demographic coding errors affecting all Black and Asian-American voters, developed geospatial ML",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Discovered systematic demographic coding errors affecting 50M voters, developed geospatial ML",
    2)
Write-Output "Summary replace: $r1"

# --- 2. Siege Analytics bullet point -----------------------------------
# First normalise the surrounding text (drop "all Black and Asian-American"),
# leaving a plain "50M" token in its place.
$r2 = $d.Content.Find.Execute(
    "Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Discovered systematic race coding errors affecting 50M voters, developed geospatial machine learning",
    2)
Write-Output "Bullet replace: $r2"

# Re-find the now-unique phrase starting with that "50M" token (the
# professional-summary paragraph also contains a "50M voters, developed
# geospatial" string, so the search text must be long enough to only
# match the bullet point's "... machine learning" wording), then shrink
# the found range down to just the 3 characters "50M" so the surrounding
# text is left in its own (unformatted) runs when we restyle it. This
# splits the run into: plain / bold+colored "50M" / plain, matching the
# existing "23%"/"64%" run pattern in this paragraph.
$metric = $d.Content
$r3 = $metric.Find.Execute("50M voters, developed geospatial machine learning")
Write-Output "Bullet metric find: $r3"
$metric.End = $metric.Start + 3
Write-Output "Bullet metric text: $($metric.Text)"
$metric.Font.Bold = 1
$metric.Font.Color = 5258796

# --- 3. "Impact:" project statement ------------------------------------
$r4 = $d.Content.Find.Execute(
    "Impact: Corrected demographic data affecting all Black and Asian-American voters, improved electoral prediction accuracy by 22%",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Impact: Corrected demographic data affecting 50M voters nationwide, improved electoral prediction accuracy by 22%",
    2)
Write-Output "Impact replace: $r4"
